$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A69").Value = 0
$ws.Range("A177").Value = 0
$ws.Range("A209").Value = 1
$ws.Range("A262").Value = 1
$ws.Range("A272").Value = 1
$ws.Range("A293").Value = 1
$ws.Range("A338").Value = 0
$ws.Range("A362").Value = 1
$ws.Range("A383").Value = 0
$ws.Range("A407").Value = 1
$ws.Range("A413").Value = 0
$ws.Range("A419").Value = 1
$ws.Range("A508").Value = 0
$ws.Range("A681").Value = 0
$ws.Range("A749").Value = 1
$ws.Range("A804").Value = 1
$ws.Range("A841").Value = 0
$ws.Range("A1253:A1268").Value = 0
$ws.Range("A1279").Value = 0
$ws.Range("A1292:A1293").Value = 0
$ws.Range("A1303:A1304").Value = 0
$ws.Range("A1309").Value = 0
$ws.Range("A1311").Value = 1
$ws.Range("A1316").Value = 1
$ws.Range("A1318").Value = 0
$ws.Range("A1320").Value = 1
$ws.Range("A1324").Value = 1
$ws.Range("A1331").Value = 1
$ws.Range("A1333").Value = 1
$ws.Range("A1334").Value = 0
$ws.Range("A1349").Value = 0
$ws.Range("A1359:A1360").Value = 1
$ws.Range("A1362").Value = 0
$ws.Range("A1364").Value = 0
$ws.Range("A1366").Value = 0
$ws.Range("A1373").Value = 0
$ws.Range("A1375").Value = 1
$ws.Range("A1382").Value = 1
$ws.Range("A1383").Value = 0
$ws.Range("A1395").Value = 0
$ws.Range("A1398").Value = 0
$ws.Range("A1400").Value = 1
$ws.Range("A1407").Value = 0
$ws.Range("A1415").Value = 1
$ws.Range("A1416").Value = 0
$ws.Range("A1423").Value = 1
$ws.Range("A1431").Value = 1
$ws.Range("A1437").Value = 0
$ws.Range("A1439").Value = 1
$ws.Range("A1441").Value = 1
$ws.Range("A1443:A1449").Value = 1
$ws.Range("A1600").Value = 1
$ws.Range("A1631").Value = 1
$ws.Range("A1634").Value = 0
$ws.Range("A1647").Value = 1
$ws.Range("A1654").Value = 0
$ws.Range("A1659").Value = 0
$ws.Range("A1663").Value = 0
$ws.Range("A1675").Value = 0
$ws.Range("A1717").Value = 1
$ws.Range("A1725").Value = 0
$ws.Range("A1753").Value = 1
$ws.Range("A1760:A1761").Value = 1
$ws.Range("A1778").Value = 0
$ws.Range("A1779").Value = 1
